# Extend the predicted-factors matrix with two more forecast columns
# ("t+15" and "t+16"): columns O and P, mirroring the existing layout
# (header row 1 = index, rows 2-13 = data) for the new columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column O / P header cells (row 1) need the same bold/centered/bordered
# style as the existing header cells (e.g. N1). Copy N1's formatting onto
# the new header cells first, then stamp in the numeric header values.
$ws.Range("N1").Copy()
$ws.Range("O1:P1").PasteSpecial(-4122)

$ws.Range("O1").Value = 14
$ws.Range("P1").Value = 15

# Data rows 2-13: plain (default-style) numeric values, same as the rest
# of the matrix.
$ws.Range("O2").Value = -0.868773331815523
$ws.Range("P2").Value = -0.6413875193420462

$ws.Range("O3").Value = -0.4714219364225961
$ws.Range("P3").Value = -0.4052797520108004

$ws.Range("O4").Value = 0.03844560962470986
$ws.Range("P4").Value = -0.008402650748879321

$ws.Range("O5").Value = 0.429479825211166
$ws.Range("P5").Value = 0.390726352879284

$ws.Range("O6").Value = -0.3125915767930236
$ws.Range("P6").Value = -0.3041813948739353

$ws.Range("O7").Value = -0.1614690575234892
$ws.Range("P7").Value = -0.1613905153886959

$ws.Range("O8").Value = -0.4577403821099732
$ws.Range("P8").Value = -0.4502764201483659

$ws.Range("O9").Value = 0.005072445241089941
$ws.Range("P9").Value = 0.004355462897618633

$ws.Range("O10").Value = 0.003569810520481088
$ws.Range("P10").Value = 0.004236458866885789

$ws.Range("O11").Value = 0.01035051486891284
$ws.Range("P11").Value = 0.009249860421837208

$ws.Range("O12").Value = -0.01941490639839046
$ws.Range("P12").Value = -0.01912168124808197

$ws.Range("O13").Value = 0.005174084191942426
$ws.Range("P13").Value = 0.004713785164572149
